$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "id_electricity_consumption" column header to "id_electricity".
# This updates the cell value, the shared string table, and the Excel Table
# column name (ListObject headers follow the header row cell text).
$ws.Range("B1").Value = "id_electricity"

# Update the active selection from O10 to B1.
$ws.Range("B1").Select()
